$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.025494267911823
$ws.Cells.Item(2, 4).Value = 1.029622293557509
$ws.Cells.Item(2, 5).Value = 1.025842403853385
$ws.Cells.Item(2, 6).Value = 1.023994663069089
$ws.Cells.Item(2, 9).Value = 1.032300846081392
$ws.Cells.Item(2, 10).Value = 1.030662975878256
$ws.Cells.Item(2, 11).Value = 1.032435670204791
$ws.Cells.Item(2, 12).Value = 1.028666791903852
$ws.Cells.Item(2, 13).Value = 1.026824466054024

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.026368869764339
$ws.Cells.Item(3, 4).Value = 1.03025259309959
$ws.Cells.Item(3, 5).Value = 1.026582943201328
$ws.Cells.Item(3, 6).Value = 1.025503485680592
$ws.Cells.Item(3, 9).Value = 1.032482231056694
$ws.Cells.Item(3, 10).Value = 1.031177321845912
$ws.Cells.Item(3, 11).Value = 1.032874632585069
$ws.Cells.Item(3, 12).Value = 1.029214901078311
$ws.Cells.Item(3, 13).Value = 1.028138375758012

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.026934971342464
$ws.Cells.Item(4, 4).Value = 1.030660565935474
$ws.Cells.Item(4, 5).Value = 1.027062656855838
$ws.Cells.Item(4, 6).Value = 1.026479985932026
$ws.Cells.Item(4, 9).Value = 1.03259848946209
$ws.Cells.Item(4, 10).Value = 1.031509686634143
$ws.Cells.Item(4, 11).Value = 1.033158108933363
$ws.Cells.Item(4, 12).Value = 1.029569435292668
$ws.Cells.Item(4, 13).Value = 1.028988266592921

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.027173002000022
$ws.Cells.Item(5, 4).Value = 1.030832107136805
$ws.Cells.Item(5, 5).Value = 1.027264455693699
$ws.Cells.Item(5, 6).Value = 1.026890556546702
$ws.Cells.Item(5, 9).Value = 1.032647098525973
$ws.Cells.Item(5, 10).Value = 1.031649304348574
$ws.Cells.Item(5, 11).Value = 1.033277147422895
$ws.Cells.Item(5, 12).Value = 1.02971845003227
$ws.Cells.Item(5, 13).Value = 1.029345493173619

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.027212970816792
$ws.Cells.Item(6, 4).Value = 1.030860911356725
$ws.Cells.Item(6, 5).Value = 1.027298346040234
$ws.Cells.Item(6, 6).Value = 1.026959496206854
$ws.Cells.Item(6, 9).Value = 1.032655244596936
$ws.Cells.Item(6, 10).Value = 1.03167274040739
$ws.Cells.Item(6, 11).Value = 1.033297126579846
$ws.Cells.Item(6, 12).Value = 1.029743468402982
$ws.Cells.Item(6, 13).Value = 1.02940546922977

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.026938151757118
$ws.Cells.Item(7, 4).Value = 1.030662857963007
$ws.Cells.Item(7, 5).Value = 1.027065352803054
$ws.Cells.Item(7, 6).Value = 1.026485471795389
$ws.Cells.Item(7, 9).Value = 1.03259914002421
$ws.Cells.Item(7, 10).Value = 1.031511552639446
$ws.Cells.Item(7, 11).Value = 1.033159700062453
$ws.Cells.Item(7, 12).Value = 1.029571426559156
$ws.Cells.Item(7, 13).Value = 1.028993040132974

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.025789807018747
$ws.Cells.Item(8, 4).Value = 1.029835279393523
$ws.Cells.Item(8, 5).Value = 1.026092561238316
$ws.Cells.Item(8, 6).Value = 1.024504539659672
$ws.Cells.Item(8, 9).Value = 1.032362375606269
$ws.Cells.Item(8, 10).Value = 1.030836894762308
$ws.Cells.Item(8, 11).Value = 1.032584135590783
$ws.Cells.Item(8, 12).Value = 1.028852054364982
$ws.Cells.Item(8, 13).Value = 1.027268571499021

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.023767644504359
$ws.Cells.Item(9, 4).Value = 1.028377993833413
$ws.Cells.Item(9, 5).Value = 1.024382516626543
$ws.Cells.Item(9, 6).Value = 1.02101516219413
$ws.Cells.Item(9, 9).Value = 1.031936678954321
$ws.Cells.Item(9, 10).Value = 1.029644619602196
$ws.Cells.Item(9, 11).Value = 1.031565637707621
$ws.Cells.Item(9, 12).Value = 1.027583463769089
$ws.Cells.Item(9, 13).Value = 1.024227410010787

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.02242048311582
$ws.Cells.Item(10, 4).Value = 1.02740720677411
$ws.Cells.Item(10, 5).Value = 1.023245321255946
$ws.Cells.Item(10, 6).Value = 1.01868950889748
$ws.Cells.Item(10, 9).Value = 1.031647190476262
$ws.Cells.Item(10, 10).Value = 1.02884747505341
$ws.Cells.Item(10, 11).Value = 1.030883789958764
$ws.Cells.Item(10, 12).Value = 1.026737118447285
$ws.Cells.Item(10, 13).Value = 1.022198117753419

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.02183737702411
$ws.Cells.Item(11, 4).Value = 1.026987031218647
$ws.Cells.Item(11, 5).Value = 1.022753585875489
$ws.Cells.Item(11, 6).Value = 1.017682558812417
$ws.Cells.Item(11, 9).Value = 1.031520492088799
$ws.Cells.Item(11, 10).Value = 1.02850176205755
$ws.Cells.Item(11, 11).Value = 1.030587872124912
$ws.Cells.Item(11, 12).Value = 1.026370500899461
$ws.Cells.Item(11, 13).Value = 1.021318923492812

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.021620819386394
$ws.Cells.Item(12, 4).Value = 1.02683098746269
$ws.Cells.Item(12, 5).Value = 1.02257103608007
$ws.Cells.Item(12, 6).Value = 1.017308538603711
$ws.Cells.Item(12, 9).Value = 1.031473228355481
$ws.Cells.Item(12, 10).Value = 1.028373267330379
$ws.Cells.Item(12, 11).Value = 1.03047785430962
$ws.Cells.Item(12, 12).Value = 1.026234301544342
$ws.Cells.Item(12, 13).Value = 1.020992272568825

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.021667270194368
$ws.Cells.Item(13, 4).Value = 1.026864458090243
$ws.Cells.Item(13, 5).Value = 1.022610188974348
$ws.Cells.Item(13, 6).Value = 1.01738876702376
$ws.Cells.Item(13, 9).Value = 1.03148337573842
$ws.Cells.Item(13, 10).Value = 1.028400833571282
$ws.Cells.Item(13, 11).Value = 1.030501458057609
$ws.Cells.Item(13, 12).Value = 1.026263517724264
$ws.Cells.Item(13, 13).Value = 1.02106234394648

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.021819475608254
$ws.Cells.Item(14, 4).Value = 1.026974132013807
$ws.Cells.Item(14, 5).Value = 1.022738494154446
$ws.Cells.Item(14, 6).Value = 1.017651642094744
$ws.Cells.Item(14, 9).Value = 1.031516589379753
$ws.Cells.Item(14, 10).Value = 1.028491142305174
$ws.Cells.Item(14, 11).Value = 1.03057878007111
$ws.Cells.Item(14, 12).Value = 1.026359243049819
$ws.Cells.Item(14, 13).Value = 1.021291924062699

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.021913258892491
$ws.Cells.Item(15, 4).Value = 1.027041709494151
$ws.Cells.Item(15, 5).Value = 1.022817560832392
$ws.Cells.Item(15, 6).Value = 1.017813608729249
$ws.Cells.Item(15, 9).Value = 1.031537026598346
$ws.Cells.Item(15, 10).Value = 1.028546773700129
$ws.Cells.Item(15, 11).Value = 1.03062640737613
$ws.Cells.Item(15, 12).Value = 1.026418219775857
$ws.Cells.Item(15, 13).Value = 1.021433365348028

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.022459186686942
$ws.Cells.Item(16, 4).Value = 1.027435096342843
$ws.Cells.Item(16, 5).Value = 1.023277970460826
$ws.Cells.Item(16, 6).Value = 1.018756337964193
$ws.Cells.Item(16, 9).Value = 1.031655570648439
$ws.Cells.Item(16, 10).Value = 1.028870407441053
$ws.Cells.Item(16, 11).Value = 1.030903414899692
$ws.Cells.Item(16, 12).Value = 1.026761446639821
$ws.Cells.Item(16, 13).Value = 1.022256456135106

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.022801692882923
$ws.Cells.Item(17, 4).Value = 1.027681906763554
$ws.Cells.Item(17, 5).Value = 1.023566955268178
$ws.Cells.Item(17, 6).Value = 1.019347703092766
$ws.Cells.Item(17, 9).Value = 1.031729569455787
$ws.Cells.Item(17, 10).Value = 1.029073268812648
$ws.Cells.Item(17, 11).Value = 1.031076994429245
$ws.Cells.Item(17, 12).Value = 1.026976705415613
$ws.Cells.Item(17, 13).Value = 1.022772623053222

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02300149251266
$ws.Cells.Item(18, 4).Value = 1.027825884575072
$ws.Cells.Item(18, 5).Value = 1.023735580553873
$ws.Cells.Item(18, 6).Value = 1.019692643766911
$ws.Cells.Item(18, 9).Value = 1.031772601635202
$ws.Cells.Item(18, 10).Value = 1.029191541901546
$ws.Cells.Item(18, 11).Value = 1.031178175446175
$ws.Cells.Item(18, 12).Value = 1.027102248253101
$ws.Cells.Item(18, 13).Value = 1.023073646789193

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.023069622670282
$ws.Cells.Item(19, 4).Value = 1.027874980224735
$ws.Cells.Item(19, 5).Value = 1.023793088476459
$ws.Cells.Item(19, 6).Value = 1.019810261032827
$ws.Cells.Item(19, 9).Value = 1.031787252427436
$ws.Cells.Item(19, 10).Value = 1.029231861053739
$ws.Cells.Item(19, 11).Value = 1.031212664523429
$ws.Cells.Item(19, 12).Value = 1.02714505273471
$ws.Cells.Item(19, 13).Value = 1.023176280145216

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.022764942977025
$ws.Cells.Item(20, 4).Value = 1.027655424528025
$ws.Cells.Item(20, 5).Value = 1.023535943169337
$ws.Cells.Item(20, 6).Value = 1.019284254481995
$ws.Cells.Item(20, 9).Value = 1.031721643536614
$ws.Cells.Item(20, 10).Value = 1.029051509137721
$ws.Cells.Item(20, 11).Value = 1.03105837770736
$ws.Cells.Item(20, 12).Value = 1.026953611627432
$ws.Cells.Item(20, 13).Value = 1.022717248191797

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.021774653968342
$ws.Cells.Item(21, 4).Value = 1.026941834992789
$ws.Cells.Item(21, 5).Value = 1.022700708637639
$ws.Cells.Item(21, 6).Value = 1.017574231816845
$ws.Cells.Item(21, 9).Value = 1.031506814373332
$ws.Cells.Item(21, 10).Value = 1.028464550882709
$ws.Cells.Item(21, 11).Value = 1.030556013443552
$ws.Cells.Item(21, 12).Value = 1.026331054911206
$ws.Cells.Item(21, 13).Value = 1.021224320656898

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.021152216416243
$ws.Cells.Item(22, 4).Value = 1.026493336240646
$ws.Cells.Item(22, 5).Value = 1.022176157652847
$ws.Cells.Item(22, 6).Value = 1.016499104223017
$ws.Cells.Item(22, 9).Value = 1.031370572344432
$ws.Cells.Item(22, 10).Value = 1.028095035345997
$ws.Cells.Item(22, 11).Value = 1.030239574340785
$ws.Cells.Item(22, 12).Value = 1.025939505916215
$ws.Cells.Item(22, 13).Value = 1.020285198477666

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.021482163611205
$ws.Cells.Item(23, 4).Value = 1.026731078148558
$ws.Cells.Item(23, 5).Value = 1.022454175474018
$ws.Cells.Item(23, 6).Value = 1.017069048421856
$ws.Cells.Item(23, 9).Value = 1.031442907694298
$ws.Cells.Item(23, 10).Value = 1.028290967150498
$ws.Cells.Item(23, 11).Value = 1.030407379804077
$ws.Cells.Item(23, 12).Value = 1.026147084953182
$ws.Cells.Item(23, 13).Value = 1.020783089901524

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.022781548610935
$ws.Cells.Item(24, 4).Value = 1.02766739065642
$ws.Cells.Item(24, 5).Value = 1.023549956000057
$ws.Cells.Item(24, 6).Value = 1.019312924154017
$ws.Cells.Item(24, 9).Value = 1.031725225320211
$ws.Cells.Item(24, 10).Value = 1.029061341560401
$ws.Cells.Item(24, 11).Value = 1.031066790003413
$ws.Cells.Item(24, 12).Value = 1.026964046759381
$ws.Cells.Item(24, 13).Value = 1.022742269854414

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.024290257077636
$ws.Cells.Item(25, 4).Value = 1.028754611257462
$ws.Cells.Item(25, 5).Value = 1.024824108352122
$ws.Cells.Item(25, 6).Value = 1.021917125289941
$ws.Cells.Item(25, 9).Value = 1.032047735924609
$ws.Cells.Item(25, 10).Value = 1.029953257143897
$ws.Cells.Item(25, 11).Value = 1.031829448066555
$ws.Cells.Item(25, 12).Value = 1.027911536033721
$ws.Cells.Item(25, 13).Value = 1.025013933357986

Write-Output "Updated vm_pu values for rows 2-25 (380 kV case)"
